# Apply 2022-06-10 data update to Fonds de solidarite Volet 1 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(13, 3).Value = 187870   # C13: nombre_aides 187867 -> 187870
$ws.Cells.Item(13, 5).Value = 1168643288   # E13: montant_total 1168405040 -> 1168643288
$ws.Cells.Item(91, 3).Value = 18887   # C91: nombre_aides 18886 -> 18887
$ws.Cells.Item(91, 5).Value = 75398735   # E91: montant_total 75390235 -> 75398735
$ws.Cells.Item(115, 3).Value = 81811   # C115: nombre_aides 81810 -> 81811
$ws.Cells.Item(115, 5).Value = 436758478   # E115: montant_total 436757089 -> 436758478
$ws.Cells.Item(121, 3).Value = 1306458   # C121: nombre_aides 1306453 -> 1306458
$ws.Cells.Item(121, 5).Value = 2275688370   # E121: montant_total 2275669602 -> 2275688370
$ws.Cells.Item(127, 3).Value = 9163   # C127: nombre_aides 9162 -> 9163
$ws.Cells.Item(127, 5).Value = 110743559   # E127: montant_total 110643825 -> 110743559
$ws.Cells.Item(129, 3).Value = 633923   # C129: nombre_aides 633912 -> 633923
$ws.Cells.Item(129, 5).Value = 3437003883   # E129: montant_total 3436652220 -> 3437003883
$ws.Cells.Item(132, 3).Value = 586118   # C132: nombre_aides 586101 -> 586118
$ws.Cells.Item(132, 5).Value = 3475155346   # E132: montant_total 3474796202 -> 3475155346
$ws.Cells.Item(134, 3).Value = 7028   # C134: nombre_aides 7027 -> 7028
$ws.Cells.Item(134, 5).Value = 16907539   # E134: montant_total 16897539 -> 16907539
$ws.Cells.Item(136, 3).Value = 26709   # C136: nombre_aides 26708 -> 26709
$ws.Cells.Item(136, 5).Value = 144431319   # E136: montant_total 144428758 -> 144431319
$ws.Cells.Item(161, 3).Value = 555   # C161: nombre_aides 554 -> 555
$ws.Cells.Item(161, 5).Value = 1818521   # E161: montant_total 1651122 -> 1818521
$ws.Cells.Item(171, 3).Value = 95830   # C171: nombre_aides 95828 -> 95830
$ws.Cells.Item(171, 5).Value = 490707985   # E171: montant_total 490704985 -> 490707985
$ws.Cells.Item(174, 3).Value = 40453   # C174: nombre_aides 40452 -> 40453
$ws.Cells.Item(174, 5).Value = 240012967   # E174: montant_total 239945952 -> 240012967
$ws.Cells.Item(186, 3).Value = 236844   # C186: nombre_aides 236842 -> 236844
$ws.Cells.Item(186, 5).Value = 1190209648   # E186: montant_total 1190190297 -> 1190209648
$ws.Cells.Item(215, 3).Value = 230266   # C215: nombre_aides 230265 -> 230266
$ws.Cells.Item(215, 5).Value = 408781748   # E215: montant_total 408775539 -> 408781748

$wb.Save()
